$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D12").Value = 4
$ws.Range("D13").Value = 1
$ws.Range("F13").Value = 30
$ws.Range("D12").Select() | Out-Null
